# Updates the "cryptos" worksheet: refreshed Price / Volume(1h) figures
# scraped from coinranking.com, plus a reordering of the ONDO /
# InjectiveProtocol rows (47-48).
#
# Note: several "Price" values look numeric (e.g. "15.00", "0.999") but
# must stay as literal text, matching the original inline strings, so a
# leading apostrophe is used to force text entry for those that would
# otherwise be auto-converted to a number by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.301.82"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.566.30"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'607.45"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'145.23"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").Value = "3.562.02"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").Value = "'7.96"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'0.412"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "4.171.83"
$ws.Range("E13").Value = "  +2.95%  "
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "'30.30"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "3.570.11"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "66.350.11"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("E19").Value = "  +9.58%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'15.00"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").Value = "'431.44"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("D24").Value = "'78.54"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "3.711.87"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D27").Value = "'0.0000121"
$ws.Range("E27").Value = "  +6.12%  "
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "'8.06"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "'9.25"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'1.50"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").Value = "3.564.52"
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("D35").Value = "'25.45"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "'7.91"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").Value = "'5.67"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'171.94"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").Value = "'0.0859"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "'5.29"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").Value = "'1.94"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'26.31"
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.22"
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("D51").Value = "'0.951"
$ws.Range("E51").Value = "  -0.90%  "
